# The commit swaps the `name` attribute of the embedded-picture docPr /
# cNvPr elements for the Pearson logo (in both footers) and the BTec logo
# (in both headers). `name` is not exposed as a writable InlineShape
# property in this COM host, so we round-trip the whole package through
# Document.WordOpenXML and patch the four `name="..."` occurrences
# directly in the OOXML text, then write it back.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

# Pearson logo (footer1.xml id="2", footer2.xml id="4"): image2.png -> image1.png
$xml = $xml.Replace('name="image2.png"', 'name="image1.png"')

# BTec logo (header1.xml id="1", header2.xml id="3"): image1.jpg -> image2.jpg
$xml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')

$d.WordOpenXML = $xml

Write-Output "renamed inline picture docPr/cNvPr names"
